# This script updates the "想去人数" (interested-count) figures and the
# status of the "抚州·第七届FZ动漫文化节" event (now cancelled) on both the
# "展览" and "全部类型" worksheets, which hold duplicate data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F ("想去人数")
$fUpdates = @{
    2  = 21
    3  = 1807
    5  = 780
    7  = 107
    12 = 148
    15 = 4209
    16 = 13
    17 = 28
    18 = 462
    19 = 398
    20 = 971
    21 = 1417
    22 = 361
    23 = 39
    24 = 43
    25 = 1952
    28 = 94
    29 = 193
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    # Row 9: event cancelled -> rename and mark ticket price as unavailable
    $ws.Range("C9").Value = "抚州·第七届FZ动漫文化节（取消)"
    $ws.Range("G9").Value = "不可售"
}
